$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = "[프로그래머스 - Python] 전화번호 목록"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%A8%B8%EC%8A%A4-Python-%EC%A0%84%ED%99%94%EB%B2%88%ED%98%B8-%EB%AA%A9%EB%A1%9D"

$ws.Range("D8").Value = "카카오브레인"

$ws.Range("D9").Value = "통계 문맹 + DNN 마니악이 되는 이유"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/stat-illiteracy-engineers-false-reasons/#utm_source=rss&utm_medium=rss&utm_campaign=stat-illiteracy-engineers-false-reasons"

$ws.Range("D36").Value = "Revisiting CNNs"
$ws.Range("E36").Value = "http://dmqm.korea.ac.kr/activity/seminar/356"

$ws.Range("D41").Value = "DevOps 에 대한 이해"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/devops/"

$ws.Range("D46").Value = "맹장염 (충수염)"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/410"
